$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add the two new worksheets (after the existing CreateNewLead sheet)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "CreatePriceList"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "CreatePriceListItems"

# ------------------------------------------------------------------
# Sheet 3: CreatePriceListItems
# ------------------------------------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws3.Columns.Item(2).ColumnWidth = 15.166666666666666

$ws3.Range("A1").Value = "Product"
$ws3.Range("B1").Value = "Unit"
$ws3.Range("C1").Value = "Amount"
$ws3.Range("A1:C1").Font.Bold = $true

$ws3.Range("A2").Value = "ArmBand 100"
$ws3.Range("B2").Value = "Basic Package"
$ws3.Range("C2").Value = 100

$ws3.Range("A3").Value = "ArmBand 150"
$ws3.Range("B3").Value = "Basic Package"
$ws3.Range("C3").Value = 150

$ws3.Range("C4").Select()

# ------------------------------------------------------------------
# Sheet 2: CreatePriceList
# ------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 10.333333333333334
$ws2.Columns.Item(2).ColumnWidth = 8.333333333333334
$ws2.Columns.Item(3).ColumnWidth = 7.833333333333333
$ws2.Columns.Item(4).ColumnWidth = 24.666666666666668

# Columns B and C hold dates that must stay as plain text ("@" format),
# matching how the source workbook stores them - set the format before
# any value is written so Excel does not auto-convert the text to a
# date serial number.
$ws2.Columns.Item(2).NumberFormat = "@"
$ws2.Columns.Item(3).NumberFormat = "@"
$ws2.Range("B1:C1").NumberFormat = "@"
$ws2.Range("B2:C2").NumberFormat = "@"

$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "StartDate"
$ws2.Range("C1").Value = "EndDate"
$ws2.Range("D1").Value = "Description"
$ws2.Range("A1:D1").Font.Bold = $true

$ws2.Range("A2").Value = "TestPriceList"
$ws2.Range("B2").Value = "2/1/2020"
$ws2.Range("C2").Value = "2/1/2030"
$ws2.Range("D2").Value = "Price list for testing purposes"

$ws2.Range("C3").Select()

# ------------------------------------------------------------------
# Sheet 1: CreateNewLead keeps its data; it simply stops being the
# tab that is selected when the workbook is (re)opened - CreatePriceList
# becomes the active sheet/tab instead.
# ------------------------------------------------------------------
$ws1.Range("F3").Select()
$ws2.Activate()
